# chore: remove invoice feature; fix monthly export weekly sheets + admin cache-bust
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Weekly Timesheet")

# --- Row 2: 2026-01-26 ---
$ws.Range("B2").Value = "Davis"
$ws.Range("C2").Value = 6.5
$ws.Range("F2").Value = 650

# --- Row 3: 2026-01-27 ---
$ws.Range("B3").Value = "Hopkins"
$ws.Range("C3").Value = 9
$ws.Range("F3").Value = 900

# --- Row 4: 2026-01-28 ---
$ws.Range("B4").Value = "Caputo (Insp.)"
$ws.Range("C4").Value = 8.5
$ws.Range("F4").Value = 850

# --- Row 5: 2026-01-29 (hours/rate/total unchanged; client name text updated) ---
$ws.Range("B5").Value = "Keevil"

# --- Row 6: 2026-01-30 ---
$ws.Range("B6").Value = "McFarland"
$ws.Range("C6").Value = 8
$ws.Range("F6").Value = 800

# --- Row 8: SUBTOTAL ---
$ws.Range("C8").Value = 38
$ws.Range("D8").Value = "Reg: 38 / OT: 0"
$ws.Range("F8").Value = 3800
